# Update the cryptos price/volume snapshot (GitHub Actions scheduled refresh).
#
# Column D ("Price") cells hold numeric-looking text (e.g. "33.873.10",
# "31.84", "0.0509") that must stay stored as TEXT, exactly like the
# original cells. Excel's COM `.Value` setter auto-detects plain numeric
# strings and converts them to real numbers (and even normalizes them,
# e.g. "13.50" -> 13.5), which would not match the source data. Prefixing
# the value with a leading apostrophe forces Excel to keep it as text
# (quote-prefixed), but that also flips the cell's style to a new
# "quote prefix" style index - so we snapshot/restore the original style
# right after the write to leave formatting untouched.
#
# Column E ("Volume(1h)") cells are already non-numeric text (padded with
# spaces, e.g. "  -0.88%  "), so a plain assignment keeps them as text
# with no style side effects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($row, $value) {
    $cell = $ws.Cells.Item($row, 4)
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = $value
}

# row -> (new Price, new Volume(1h)); omit the side when unchanged.
Set-TextPrice 2 "33.873.10"
Set-Volume    2 "  -0.88%  "

Set-TextPrice 3 "1.780.28"
Set-Volume    3 "  -1.27%  "

Set-Volume    4 "  +0.05%  "

Set-TextPrice 5 "224.08"
Set-Volume    5 "  +0.30%  "

Set-Volume    6 "  -1.21%  "

Set-Volume    7 "  +0.06%  "

Set-TextPrice 8 "31.84"
Set-Volume    8 "  -1.45%  "

Set-Volume    9 "  +0.71%  "

Set-Volume    10 "  -5.79%  "

Set-Volume    11 "  +0.90%  "

Set-Volume    12 "  -1.27%  "

Set-TextPrice 13 "11.21"
Set-Volume    13 "  +1.80%  "

Set-TextPrice 14 "1.778.55"
Set-Volume    14 "  -0.40%  "

Set-TextPrice 15 "33.876.83"
Set-Volume    15 "  -0.95%  "

Set-Volume    16 "  -3.59%  "

Set-Volume    17 "  -2.27%  "

Set-TextPrice 18 "66.72"
Set-Volume    18 "  -2.69%  "

Set-TextPrice 19 "238.88"
Set-Volume    19 "  -3.71%  "

Set-Volume    20 "  -2.21%  "

Set-Volume    21 "  +0.08%  "

Set-TextPrice 22 "10.59"
Set-Volume    22 "  -3.40%  "

Set-Volume    23 "  -2.55%  "

Set-Volume    24 "  -2.57%  "

Set-TextPrice 25 "160.72"
Set-Volume    25 "  +0.82%  "

Set-Volume    26 "  -1.12%  "

Set-Volume    27 "  -3.41%  "

Set-Volume    29 "  +0.17%  "

Set-Volume    30 "  +0.58%  "

Set-TextPrice 31 "0.0509"
Set-Volume    31 "  -3.05%  "

Set-Volume    32 "  -3.86%  "

Set-Volume    33 "  +0.05%  "

Set-Volume    34 "  -2.15%  "

Set-TextPrice 35 "1.390.73"
Set-Volume    35 "  -1.94%  "

Set-Volume    36 "  -2.93%  "

Set-Volume    37 "  -1.57%  "

Set-Volume    38 "  -1.11%  "

Set-Volume    39 "  +4.68%  "

Set-Volume    40 "  +0.07%  "

Set-Volume    41 "  -3.63%  "

Set-TextPrice 42 "78.38"
Set-Volume    42 "  -2.74%  "

Set-TextPrice 43 "13.50"
Set-Volume    43 "  +12.23%  "

Set-Volume    44 "  -3.20%  "

Set-Volume    45 "  +11.50%  "

Set-TextPrice 46 "0.0509"
Set-Volume    46 "  +2.64%  "

Set-Volume    47 "  +2.95%  "

Set-Volume    48 "  -1.60%  "

Set-TextPrice 49 "106.35"
Set-Volume    49 "  -1.39%  "

Set-TextPrice 50 "1.937.81"
Set-Volume    50 "  -1.37%  "

Set-Volume    51 "  +0.08%  "
